# Applies the "Add files via upload" edit: appends 15 new rows (375-389) of
# PED endpoint records to Sheet2, introducing 2 brand-new cell styles along
# the way, and updates the sheet's view/pane/selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Step 1: stamp formatting on every new cell FIRST (via copy / paste-special
# "formats only"), using existing cells elsewhere in the workbook as style
# donors. This reproduces the exact cellXfs entries (including the two new
# ones) without disturbing the shared-strings pool, whose append order we
# control precisely in step 2 below.
# ---------------------------------------------------------------------------

# Donor for style "9"  (numFmt 49 / font 0 / no border / center) -> A,B,C cols
$donor9 = $ws2.Range("A2")
# Donor for style "11" (numFmt 49 / font 1 / no border / center) -> D,F,G cols
$donor11 = $ws2.Range("F2")
# Donor for the workbook default style "0" (General, no border, no align)
$donor0 = $ws1.Range("A3")
# Donor for style "12" (numFmt 49 / font 0 / border 1 / center) -> A,B,C (row389)
$donor12 = $ws2.Range("A299")
# Donor for style "13" (numFmt 49 / font 1 / border 1 / center) -> D,F,G (row389)
$donor13 = $ws2.Range("F299")
# Donor for style "18" (numFmt 49 / font 1 / border 1 / no align) -> H,I,J (row389)
$donor18 = $ws2.Range("I299")
# Donor used as the base for the two brand-new styles (General, font 0)
$donorGeneralNoBorder = $ws1.Range("F1")
$donorGeneralBorder   = $ws1.Range("C13")

# Rows 375-377 use A/B/C -> style 9, D/F/G -> style 11, E -> default style 0
$rowsPlain = 375,376,377
foreach ($r in $rowsPlain) {
    $donor9.Copy()
    $ws2.Range("A$r:C$r").PasteSpecial(-4122)
    $donor11.Copy()
    $ws2.Range("D$r").PasteSpecial(-4122)
    $ws2.Range("F$r:G$r").PasteSpecial(-4122)
    $donor0.Copy()
    $ws2.Range("E$r").PasteSpecial(-4122)
}

# Rows 378-388 use A/B/C -> style 9, D/F/G -> style 11, E -> NEW style (left,
# general, no border) -- created once, off-screen, then reused by copy.
$donor9.Copy()
$ws2.Range("A378:C388").PasteSpecial(-4122)
$donor11.Copy()
$ws2.Range("D378:D388").PasteSpecial(-4122)
$ws2.Range("F378:G388").PasteSpecial(-4122)

$donorGeneralNoBorder.Copy()
$ws2.Range("E378").PasteSpecial(-4122)
$ws2.Range("E378").HorizontalAlignment = -4131
$ws2.Range("E378").Copy()
$ws2.Range("E379:E388").PasteSpecial(-4122)

# Row 389 (the bottom/bold row): A/B/C -> style 12, D/F/G -> style 13,
# E -> NEW style (left, general, border), H/I/J -> style 18 (blank cells).
$donor12.Copy()
$ws2.Range("A389:C389").PasteSpecial(-4122)
$donor13.Copy()
$ws2.Range("D389").PasteSpecial(-4122)
$ws2.Range("F389:G389").PasteSpecial(-4122)
$donor18.Copy()
$ws2.Range("H389:J389").PasteSpecial(-4122)

$donorGeneralBorder.Copy()
$ws2.Range("E389").PasteSpecial(-4122)
$ws2.Range("E389").HorizontalAlignment = -4131

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 2: write the cell VALUES in precisely the order the original author's
# shared-strings table was built (new unique strings are appended to the
# pool in first-use order), so xl/sharedStrings.xml matches byte-for-byte.
# ---------------------------------------------------------------------------

$ws2.Range("E375").Value = "ศูนย์พักคอย อำเภอเมือง"
$ws2.Range("E376").Value = "หอประชุม อบต.กง"
$ws2.Range("E377").Value = "อบต.ไกรนอก"
$ws2.Range("E378").Value = "อบต.ไกรใน"
$ws2.Range("E379").Value = "อบต.ไกรกลาง"
$ws2.Range("E380").Value = "ศาลาวัดบ้านใหม่สุขเกษม หลังที่ 1,2"
$ws2.Range("E381").Value = "อาคารวัฒนธรรมเทศบาลตำบลกงไกรลาศ"
$ws2.Range("E382").Value = "โรงเรียนข่อยสองนาง"
$ws2.Range("E383").Value = "อาคารศูนย์พัฒนาเด็กเล็ก อบต.ดงเดือย"
$ws2.Range("E384").Value = "วัดโบสถ์ ม.1 ต.ไกรนอก"
$ws2.Range("E385").Value = "ตลาดชุมชนบ้านโป่งแค หมู่ที่ 7 ต.ไกรนอก"
$ws2.Range("E386").Value = "อาคารเอนกประสงค์ที่สาธารณะประโยชน์ ต.กกแรต หมู่ที่ 10"
$ws2.Range("E387").Value = "อาคารผู้สูงอายุ ม.3 ต.หนองตูม"
$ws2.Range("E388").Value = "วัดหนองตูม"
$ws2.Range("E389").Value = "อาคารเอนกประสงค์เอกชน ม.8 ต.ท่าฉนวน"

$ws2.Range("A376").Value = "PED_03_64001"
$ws2.Range("A377").Value = "PED_03_64002"
$ws2.Range("A378").Value = "PED_03_64003"
$ws2.Range("A379").Value = "PED_03_64004"
$ws2.Range("A380").Value = "PED_03_64005"
$ws2.Range("A381").Value = "PED_03_64006"
$ws2.Range("A382").Value = "PED_03_64007"
$ws2.Range("A383").Value = "PED_03_64008"
$ws2.Range("A384").Value = "PED_03_64009"
$ws2.Range("A385").Value = "PED_03_64010"
$ws2.Range("A386").Value = "PED_03_64011"
$ws2.Range("A387").Value = "PED_03_64012"
$ws2.Range("A388").Value = "PED_03_64013"
$ws2.Range("A389").Value = "PED_03_64014"
$ws2.Range("A375").Value = "PED_03_62001"

$ws2.Range("D375").Value = "62001"
$ws2.Range("D377").Value = "64002"
$ws2.Range("D378").Value = "64003"
$ws2.Range("D379").Value = "64004"
$ws2.Range("D380").Value = "64005"
$ws2.Range("D381").Value = "64006"
$ws2.Range("D382").Value = "64007"
$ws2.Range("D383").Value = "64008"
$ws2.Range("D384").Value = "64009"
$ws2.Range("D385").Value = "64010"
$ws2.Range("D386").Value = "64011"
$ws2.Range("D387").Value = "64012"
$ws2.Range("D388").Value = "64013"
$ws2.Range("D389").Value = "64014"
$ws2.Range("D376").Value = "64001"

# B / C columns reuse the already-existing "PED" / "03" shared strings.
$ws2.Range("B375:B389").Value = "PED"
$ws2.Range("C375:C389").Value = "03"

# F / G columns are plain numbers.
$ws2.Range("F375").Value = 62
$ws2.Range("G375").Value = 6201
$ws2.Range("F376:F389").Value = 64
$ws2.Range("G376:G389").Value = 6404

# ---------------------------------------------------------------------------
# Step 3: sheet view housekeeping to match the post-edit state.
# ---------------------------------------------------------------------------

$ws2.Activate()
$window = $excel.ActiveWindow
$window.Zoom = 100
$window.FreezePanes = $false
$ws2.Range("A2").Select()
$window.FreezePanes = $true
$ws2.Range("C387").Select()
